$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 111, shifting existing rows 111:135 down to 112:136
$ws.Rows(111).Insert()

# Populate the newly inserted row 111 with the new weekly record
$ws.Cells.Item(111, 1).Value = 10
$ws.Cells.Item(111, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(111, 3).Value = "La Araucanía"
$ws.Cells.Item(111, 4).Value = 45275
$ws.Cells.Item(111, 5).Value = 9
$ws.Cells.Item(111, 6).Value = "Fruta"
$ws.Cells.Item(111, 7).Value = 100108
$ws.Cells.Item(111, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(111, 9).Value = 100108007
$ws.Cells.Item(111, 10).Value = "Coco"
$ws.Cells.Item(111, 11).Value = "Sin especificar"
$ws.Cells.Item(111, 12).Value = "Primera"
$ws.Cells.Item(111, 13).Value = 15
$ws.Cells.Item(111, 14).Value = 32000
$ws.Cells.Item(111, 15).Value = 32000
$ws.Cells.Item(111, 16).Value = 32000
$ws.Cells.Item(111, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(111, 18).Value = "Perú"
$ws.Cells.Item(111, 19).Value = 1600
$ws.Cells.Item(111, 20).Value = 20
